$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 151.66667
$ws.Range("I12").Value = 158.125
$ws.Range("K12").Value = 158.125
$ws.Range("M12").Value = 11.875
$ws.Range("H58").Value = 2273.375
$ws.Range("I58").Value = 1023.3333
$ws.Range("J58").Value = 3023.4
$ws.Range("K58").Value = 3069.9999
$ws.Range("L58").Value = 9070.200000000001
$ws.Range("M58").Value = -2919.9999
$ws.Range("N58").Value = -9370.200000000001
$ws.Range("H99").Value = 168.5
$ws.Range("I99").Value = 168.5
$ws.Range("K99").Value = 505.5
$ws.Range("M99").Value = 992.5
$ws.Range("H106").Value = 5084.357
$ws.Range("I106").Value = 5298.5386
$ws.Range("K106").Value = 5298.5386
$ws.Range("M106").Value = -4667.5386
$ws.Range("H118").Value = 1582
$ws.Range("I118").Value = 268.75
$ws.Range("K118").Value = 806.25
$ws.Range("M118").Value = 850.75
$ws.Range("H121").Value = 4996
$ws.Range("J121").Value = 4996
$ws.Range("L121").Value = 14988
$ws.Range("N121").Value = -18482
$ws.Range("H138").Value = 4900
$ws.Range("I138").Value = 4521.3
$ws.Range("J138").Value = 5170.5
$ws.Range("K138").Value = 13563.9
$ws.Range("L138").Value = 15511.5
$ws.Range("M138").Value = -8423.900000000001
$ws.Range("N138").Value = -25791.5

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2300
$ws.Range("I45").Value = 2375
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 2375
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1998
$ws.Range("N45").Value = -2754
$ws.Range("H63").Value = 4463.1665
$ws.Range("I63").Value = 2149.75
$ws.Range("K63").Value = 2149.75
$ws.Range("M63").Value = -1463.75
$ws.Range("H66").Value = 4463.1665
$ws.Range("I66").Value = 2149.75
$ws.Range("K66").Value = 10748.75
$ws.Range("M66").Value = -7316.75
$ws.Range("H97").Value = 1111.9375
$ws.Range("I97").Value = 972.73334
$ws.Range("K97").Value = 972.73334
$ws.Range("M97").Value = -476.73334
$ws.Range("H122").Value = 2339.6072
$ws.Range("I122").Value = 1932.7273
$ws.Range("K122").Value = 5798.1819
$ws.Range("M122").Value = -3348.1819

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1300
$ws.Range("I86").Value = 1300
$ws.Range("K86").Value = 1300
$ws.Range("M86").Value = -177
$ws.Range("H89").Value = 1300
$ws.Range("I89").Value = 1300
$ws.Range("K89").Value = 6500
$ws.Range("M89").Value = -884
$ws.Range("H94").Value = 6229
$ws.Range("I94").Value = 5874.8
$ws.Range("K94").Value = 5874.8
$ws.Range("M94").Value = -5423.8
$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2651.36
$ws.Range("J31").Value = 3368.5
$ws.Range("L31").Value = 3368.5
$ws.Range("N31").Value = -3958.5
$ws.Range("H34").Value = 2651.36
$ws.Range("J34").Value = 3368.5
$ws.Range("L34").Value = 3368.5
$ws.Range("N34").Value = -3772.5
$ws.Range("H122").Value = 2209.5715
$ws.Range("I122").Value = 2065.6775
$ws.Range("J122").Value = 3324.75
$ws.Range("K122").Value = 6197.032499999999
$ws.Range("L122").Value = 9974.25
$ws.Range("M122").Value = -3747.032499999999
$ws.Range("N122").Value = -14874.25
$ws.Range("H134").Value = 1725.6522
$ws.Range("I134").Value = 1247.1578
$ws.Range("K134").Value = 3741.4734
$ws.Range("M134").Value = -1206.4734

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 114.181816
$ws.Range("J2").Value = 57.666668
$ws.Range("L2").Value = 346.000008
$ws.Range("N2").Value = -572.000008
$ws.Range("H4").Value = 8574136
$ws.Range("I4").Value = 25715306
$ws.Range("J4").Value = 3550.2856
$ws.Range("K4").Value = 77145918
$ws.Range("L4").Value = 10650.8568
$ws.Range("M4").Value = -77145806
$ws.Range("N4").Value = -10874.8568
$ws.Range("H12").Value = 546.2105
$ws.Range("I12").Value = 1042
$ws.Range("J12").Value = 317.3846
$ws.Range("K12").Value = 3126
$ws.Range("L12").Value = 952.1537999999999
$ws.Range("M12").Value = -2953
$ws.Range("N12").Value = -1298.1538
$ws.Range("H23").Value = 162.66667
$ws.Range("I23").Value = 88.57143000000001
$ws.Range("J23").Value = 209.81818
$ws.Range("K23").Value = 265.71429
$ws.Range("L23").Value = 629.4545400000001
$ws.Range("M23").Value = -30.71429000000001
$ws.Range("N23").Value = -1099.45454
$ws.Range("H34").Value = 1560.25
$ws.Range("J34").Value = 1894.9231
$ws.Range("L34").Value = 5684.7693
$ws.Range("N34").Value = -5852.7693
$ws.Range("H39").Value = 124580.91
$ws.Range("I39").Value = 182816.36
$ws.Range("J39").Value = 66345.45
$ws.Range("K39").Value = 548449.08
$ws.Range("L39").Value = 199036.35
$ws.Range("M39").Value = -548155.08
$ws.Range("N39").Value = -199624.35
$ws.Range("H70").Value = 3497.5
$ws.Range("I70").Value = 3497.5
$ws.Range("K70").Value = 10492.5
$ws.Range("M70").Value = -10177.5
$ws.Range("H73").Value = 3497.5
$ws.Range("I73").Value = 3497.5
$ws.Range("K73").Value = 10492.5
$ws.Range("M73").Value = -9400.5
$ws.Range("H122").Value = 1467692.9
$ws.Range("I122").Value = 8064740.5
$ws.Range("J122").Value = 1682.2222
$ws.Range("K122").Value = 72582664.5
$ws.Range("L122").Value = 15139.9998
$ws.Range("M122").Value = -72580214.5
$ws.Range("N122").Value = -20039.9998
$ws.Range("H131").Value = 4782.533
$ws.Range("I131").Value = 789.25
$ws.Range("K131").Value = 2367.75
$ws.Range("M131").Value = 2672.25
$ws.Range("H140").Value = 4133.3335
$ws.Range("I140").Value = 4133.3335
$ws.Range("K140").Value = 12400.0005
$ws.Range("M140").Value = -7220.000499999998

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 32350.064
$ws.Range("I2").Value = 55607.168
$ws.Range("K2").Value = 55607.168
$ws.Range("M2").Value = -55494.168
$ws.Range("H97").Value = 1029.75
$ws.Range("I97").Value = 955.4167
$ws.Range("J97").Value = 1252.75
$ws.Range("K97").Value = 955.4167
$ws.Range("L97").Value = 1252.75
$ws.Range("M97").Value = -459.4167
$ws.Range("N97").Value = -2244.75
$ws.Range("H107").Value = 623.5
$ws.Range("I107").Value = 597
$ws.Range("K107").Value = 597
$ws.Range("M107").Value = 1323
$ws.Range("H123").Value = 85551
$ws.Range("J123").Value = 85551
$ws.Range("L123").Value = 85551
$ws.Range("N123").Value = -90451

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3224.75
$ws.Range("I7").Value = 1452
$ws.Range("K7").Value = 1452
$ws.Range("M7").Value = -1340
$ws.Range("H40").Value = 34333.332
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3272
$ws.Range("H68").Value = 9753
$ws.Range("I68").Value = 8691
$ws.Range("K68").Value = 8691
$ws.Range("M68").Value = -7942
$ws.Range("H71").Value = 9753
$ws.Range("I71").Value = 8691
$ws.Range("K71").Value = 43455
$ws.Range("M71").Value = -39711
$ws.Range("H82").Value = 2098.0715
$ws.Range("I82").Value = 2107.5
$ws.Range("K82").Value = 2107.5
$ws.Range("M82").Value = -1746.5
$ws.Range("H85").Value = 2098.0715
$ws.Range("I85").Value = 2107.5
$ws.Range("K85").Value = 2107.5
$ws.Range("M85").Value = -859.5
$ws.Range("H126").Value = 3224.75
$ws.Range("I126").Value = 1452
$ws.Range("K126").Value = 4356
$ws.Range("M126").Value = -1886

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 602.5
$ws.Range("I100").Value = 444.2
$ws.Range("K100").Value = 888.4
$ws.Range("M100").Value = -347.4
$ws.Range("H126").Value = 1812.5238
$ws.Range("I126").Value = 1418.5714
$ws.Range("K126").Value = 4255.7142
$ws.Range("M126").Value = -1785.7142
$ws.Range("H132").Value = 4399.5
$ws.Range("I132").Value = 4065.6667
$ws.Range("J132").Value = 4733.3335
$ws.Range("K132").Value = 12197.0001
$ws.Range("L132").Value = 14200.0005
$ws.Range("M132").Value = -9667.000100000001
$ws.Range("N132").Value = -19260.0005
